$d = $word.ActiveDocument
$apos = [char]0x2019

# ------------------------------------------------------------------
# "Il consulente non riesce a comunicare col server e non riesce a
#  vedere la pagina." ->
# "Il consulente non riesce a visualizzare la pagina <i>contiene
#  l'elenco degli ordini effettuati dai clienti</i>."
# ------------------------------------------------------------------

# Step 1: swap the failing-description clause for the new wording,
# leaving a unique marker where the italic clause must go.
$r1 = $d.Content.Duplicate
$r1.Find.Execute("comunicare col server e non riesce a vedere la pagina.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "visualizzare la pagina @@MARKER@@.", 2)

# Step 2: replace the marker with the italicised clause.
$r2 = $d.Content.Duplicate
$r2.Find.ClearFormatting()
$r2.Find.Replacement.ClearFormatting()
$r2.Find.Replacement.Font.Italic = $true
$r2.Find.Execute("@@MARKER@@", $false, $false, $false, $false, $false, $true, 1, $false, `
    ("contiene l" + $apos + "elenco degli ordini effettuati dai clienti"), 2)
